# Updated symbol list on Thu Dec 15 21:35:43 UTC 2022 with GitHub Actions
#
# Applies the price / row updates described by the diff against cryptos.xlsx.
# Column D holds numeric-looking values that are stored as TEXT (inline
# strings) in the workbook, so each assignment is given a leading apostrophe
# to force Excel to keep it as text (quote-prefixed) instead of coercing it
# to a real number. Columns B/C/E are plain text and are assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) corrections -----------------------------------------
$ws.Range("D2").Value  = "'259.71"
$ws.Range("D4").Value  = "'6.183"
$ws.Range("D5").Value  = "'0.06098"
$ws.Range("D6").Value  = "'6.729"
$ws.Range("D7").Value  = "'3.484"
$ws.Range("D8").Value  = "'1.363"
$ws.Range("D9").Value  = "'0.7990"
$ws.Range("D10").Value = "'0.1586"
$ws.Range("D11").Value = "'0.08037"
$ws.Range("D12").Value = "'0.03321"
$ws.Range("D13").Value = "'0.03049"
$ws.Range("D14").Value = "'0.09308"
$ws.Range("D15").Value = "'3.896"
$ws.Range("D16").Value = "'0.001696"
$ws.Range("D17").Value = "'0.04835"
$ws.Range("D19").Value = "'0.006206"
$ws.Range("D21").Value = "'0.003392"
$ws.Range("D23").Value = "'3.690"
$ws.Range("D24").Value = "'2.261"
$ws.Range("D40").Value = "'0.04586"

# --- Row 41: BKEXToken -> KickToken -----------------------------------------
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.007106"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("D42").Value = "'0.003903"

# --- Row 43: KickToken -> BKEXToken -----------------------------------------
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1118"
$ws.Range("E43").Value = "42BKEXTokenBKK"

$ws.Range("D44").Value = "'0.01065"
$ws.Range("D46").Value = "'0.00005930"
$ws.Range("D49").Value = "'0.06509"
